$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A4 label text (1st run gets ", haar" suffix)
$ws.Range('A4').Value = 'Multiresolution Histograms (1st Run - every SVM, haar)'

# Row 34 becomes the label for the 3rd run (old '2nd Run' text + ', haar')
$ws.Range('A34').Value = 'Multiresolution Histograms (2nd Run - some SVM, haar)'
$ws.Range('A34').Font.Bold = $true

# ---- Block 3: Multiresolution Histograms 3rd run (haar) rows 35-47 ----
$ws.Range('A35').Font.Bold = $true
$ws.Range('A35').Value = ''
$ws.Range('B35').Value = 'Guess'
$ws.Range('B35').Font.Bold = $true

$ws.Range('A36').Font.Bold = $true
$ws.Range('A36').Value = ''
$ws.Range('B36').Value = 'Airport'
$ws.Range('C36').Value = 'Auditorium'
$ws.Range('D36').Value = 'Bamboo'
$ws.Range('E36').Value = 'Campus'
$ws.Range('F36').Value = 'Desert'
$ws.Range('G36').Value = 'Football Field'
$ws.Range('H36').Value = 'Kitchen'
$ws.Range('I36').Value = 'Sky'
$ws.Range('J36').Value = 'Percentage'
$ws.Range('K36').Value = 'False Neg'
$ws.Range('L36').Value = 'True Pos'

$ws.Range('A37').Value = 'Airport'
$ws.Range('B37').Value = 0
$ws.Range('C37').Value = 0
$ws.Range('D37').Value = 0
$ws.Range('E37').Value = 0
$ws.Range('F37').Value = 16
$ws.Range('G37').Value = 1
$ws.Range('H37').Value = 0
$ws.Range('I37').Value = 3
$ws.Range('J37').Formula = '=SUM(B37:I37)/SUM($B$37:$I$44)'
$ws.Range('K37').Formula = '=(SUM(B37:I37) - B37) / SUM(B37:I37)'
$ws.Range('L37').Formula = '=1-K37'

$ws.Range('A38').Value = 'Auditorium'
$ws.Range('B38').Value = 0
$ws.Range('C38').Value = 2
$ws.Range('D38').Value = 1
$ws.Range('E38').Value = 1
$ws.Range('F38').Value = 1
$ws.Range('G38').Value = 8
$ws.Range('H38').Value = 4
$ws.Range('I38').Value = 3
$ws.Range('J38').Formula = '=SUM(B38:I38)/SUM($B$37:$I$44)'
$ws.Range('K38').Formula = '=(SUM(B38:I38) - C38) / SUM(B38:I38)'
$ws.Range('L38').Formula = '=1-K38'

$ws.Range('A39').Value = 'Bamboo'
$ws.Range('B39').Value = 0
$ws.Range('C39').Value = 0
$ws.Range('D39').Value = 14
$ws.Range('E39').Value = 2
$ws.Range('F39').Value = 0
$ws.Range('G39').Value = 1
$ws.Range('H39').Value = 3
$ws.Range('I39').Value = 0
$ws.Range('J39').Formula = '=SUM(B39:I39)/SUM($B$37:$I$44)'
$ws.Range('K39').Formula = '=(SUM(B39:I39) - D39) / SUM(B39:I39)'
$ws.Range('L39').Formula = '=1-K39'

$ws.Range('A40').Value = 'Campus'
$ws.Range('B40').Value = 0
$ws.Range('C40').Value = 1
$ws.Range('D40').Value = 3
$ws.Range('E40').Value = 10
$ws.Range('F40').Value = 2
$ws.Range('G40').Value = 4
$ws.Range('H40').Value = 0
$ws.Range('I40').Value = 0
$ws.Range('J40').Formula = '=SUM(B40:I40)/SUM($B$37:$I$44)'
$ws.Range('K40').Formula = '=(SUM(B40:I40) - E40) / SUM(B40:I40)'
$ws.Range('L40').Formula = '=1-K40'

$ws.Range('A41').Value = 'Desert'
$ws.Range('B41').Value = 0
$ws.Range('C41').Value = 0
$ws.Range('D41').Value = 0
$ws.Range('E41').Value = 0
$ws.Range('F41').Value = 16
$ws.Range('G41').Value = 3
$ws.Range('H41').Value = 0
$ws.Range('I41').Value = 1
$ws.Range('J41').Formula = '=SUM(B41:I41)/SUM($B$37:$I$44)'
$ws.Range('K41').Formula = '=(SUM(B41:I41) - F41) / SUM(B41:I41)'
$ws.Range('L41').Formula = '=1-K41'

$ws.Range('A42').Value = 'Football Field'
$ws.Range('B42').Value = 0
$ws.Range('C42').Value = 1
$ws.Range('D42').Value = 2
$ws.Range('E42').Value = 1
$ws.Range('F42').Value = 11
$ws.Range('G42').Value = 4
$ws.Range('H42').Value = 0
$ws.Range('I42').Value = 1
$ws.Range('J42').Formula = '=SUM(B42:I42)/SUM($B$37:$I$44)'
$ws.Range('K42').Formula = '=(SUM(B42:I42) - G42) / SUM(B42:I42)'
$ws.Range('L42').Formula = '=1-K42'

$ws.Range('A43').Value = 'Kitchen'
$ws.Range('B43').Value = 0
$ws.Range('C43').Value = 2
$ws.Range('D43').Value = 2
$ws.Range('E43').Value = 0
$ws.Range('F43').Value = 2
$ws.Range('G43').Value = 3
$ws.Range('H43').Value = 11
$ws.Range('I43').Value = 0
$ws.Range('J43').Formula = '=SUM(B43:I43)/SUM($B$37:$I$44)'
$ws.Range('K43').Formula = '=(SUM(B43:I43) - H43) / SUM(B43:I43)'
$ws.Range('L43').Formula = '=1-K43'

$ws.Range('A44').Value = 'Sky'
$ws.Range('B44').Value = 0
$ws.Range('C44').Value = 1
$ws.Range('D44').Value = 0
$ws.Range('E44').Value = 0
$ws.Range('F44').Value = 6
$ws.Range('G44').Value = 1
$ws.Range('H44').Value = 2
$ws.Range('I44').Value = 10
$ws.Range('J44').Formula = '=SUM(B44:I44)/SUM($B$37:$I$44)'
$ws.Range('K44').Formula = '=(SUM(B44:I44) - I44) / SUM(B44:I44)'
$ws.Range('L44').Formula = '=1-K44'

$ws.Range('A45').Value = 'Percentage'
$ws.Range('B45').Formula = '=SUM(B37:B44) / SUM($B$37:$I$44)'
$ws.Range('C45').Formula = '=SUM(C37:C44) / SUM($B$22:$I$29)'
$ws.Range('D45').Formula = '=SUM(D37:D44) / SUM($B$22:$I$29)'
$ws.Range('E45').Formula = '=SUM(E37:E44) / SUM($B$22:$I$29)'
$ws.Range('F45').Formula = '=SUM(F37:F44) / SUM($B$22:$I$29)'
$ws.Range('G45').Formula = '=SUM(G37:G44) / SUM($B$22:$I$29)'
$ws.Range('H45').Formula = '=SUM(H37:H44) / SUM($B$22:$I$29)'
$ws.Range('I45').Formula = '=SUM(I37:I44) / SUM($B$22:$I$29)'

$ws.Range('A46').Value = 'False Pos'
$ws.Range('C46').Formula = '=(SUM(C37:C44) - C38) / SUM(C37:C44)'
$ws.Range('D46').Formula = '=(SUM(D37:D44) - D39) / SUM(D37:D44)'
$ws.Range('E46').Formula = '=(SUM(E37:E44) - E40) / SUM(E37:E44)'
$ws.Range('F46').Formula = '=(SUM(F37:F44) - F41) / SUM(F37:F44)'
$ws.Range('G46').Formula = '=(SUM(G37:G44) - G42) / SUM(G37:G44)'
$ws.Range('H46').Formula = '=(SUM(H37:H44) - H43) / SUM(H37:H44)'
$ws.Range('I46').Formula = '=(SUM(I37:I44) - I44) / SUM(I37:I44)'
$ws.Range('K46').Value = 'Accuracy'
$ws.Range('K46').Font.Bold = $true
$ws.Range('L46').Formula = '=(B37+C38+D39+E40+F41+G42+H43+I44) / SUM(B37:I44)'

$ws.Range('A47').Value = 'True Neg'
$ws.Range('C47').Formula = '=1-C46'
$ws.Range('D47').Formula = '=1-D46'
$ws.Range('E47').Formula = '=1-E46'
$ws.Range('F47').Formula = '=1-F46'
$ws.Range('G47').Formula = '=1-G46'
$ws.Range('H47').Formula = '=1-H46'
$ws.Range('I47').Formula = '=1-I46'
$ws.Range('K47').Value = 'Total'
$ws.Range('K47').Font.Bold = $true
$ws.Range('L47').Formula = '=SUM(B37:I44)'

# ---- Block 4: Multiresolution Histograms 4th run (no SVM, bayes, haar) rows 49-62 ----
$ws.Range('A49').Value = 'Multiresolution Histograms (3nd Run - no SVM, bayes NB_THRESH = .5, haar)'
$ws.Range('A49').Font.Bold = $true

$ws.Range('B50').Value = 'Guess'
$ws.Range('B50').Font.Bold = $true

$ws.Range('B51').Value = 'Airport'
$ws.Range('C51').Value = 'Auditorium'
$ws.Range('D51').Value = 'Bamboo'
$ws.Range('E51').Value = 'Campus'
$ws.Range('F51').Value = 'Desert'
$ws.Range('G51').Value = 'Football Field'
$ws.Range('H51').Value = 'Kitchen'
$ws.Range('I51').Value = 'Sky'
$ws.Range('J51').Value = 'Percentage'
$ws.Range('K51').Value = 'False Neg'
$ws.Range('L51').Value = 'True Pos'

$ws.Range('A52').Value = 'Airport'
$ws.Range('B52').Value = 12
$ws.Range('C52').Value = 1
$ws.Range('D52').Value = 0
$ws.Range('E52').Value = 0
$ws.Range('F52').Value = 3
$ws.Range('G52').Value = 1
$ws.Range('H52').Value = 0
$ws.Range('I52').Value = 3
$ws.Range('J52').Formula = '=SUM(B52:I52)/SUM($B$22:$I$29)'
$ws.Range('K52').Formula = '=(SUM(B52:I52) - B52) / SUM(B52:I52)'
$ws.Range('L52').Formula = '=1-K52'

$ws.Range('A53').Value = 'Auditorium'
$ws.Range('B53').Value = 1
$ws.Range('C53').Value = 5
$ws.Range('D53').Value = 1
$ws.Range('E53').Value = 0
$ws.Range('F53').Value = 4
$ws.Range('G53').Value = 1
$ws.Range('H53').Value = 7
$ws.Range('I53').Value = 1
$ws.Range('J53').Formula = '=SUM(B53:I53)/SUM($B$22:$I$29)'
$ws.Range('K53').Formula = '=(SUM(B53:I53) - C53) / SUM(B53:I53)'
$ws.Range('L53').Formula = '=1-K53'

$ws.Range('A54').Value = 'Bamboo'
$ws.Range('B54').Value = 1
$ws.Range('C54').Value = 1
$ws.Range('D54').Value = 14
$ws.Range('E54').Value = 0
$ws.Range('F54').Value = 2
$ws.Range('G54').Value = 0
$ws.Range('H54').Value = 1
$ws.Range('I54').Value = 1
$ws.Range('J54').Formula = '=SUM(B54:I54)/SUM($B$22:$I$29)'
$ws.Range('K54').Formula = '=(SUM(B54:I54) - D54) / SUM(B54:I54)'
$ws.Range('L54').Formula = '=1-K54'

$ws.Range('A55').Value = 'Campus'
$ws.Range('B55').Value = 8
$ws.Range('C55').Value = 2
$ws.Range('D55').Value = 3
$ws.Range('E55').Value = 0
$ws.Range('F55').Value = 6
$ws.Range('G55').Value = 0
$ws.Range('H55').Value = 1
$ws.Range('I55').Value = 0
$ws.Range('J55').Formula = '=SUM(B55:I55)/SUM($B$22:$I$29)'
$ws.Range('K55').Formula = '=(SUM(B55:I55) - E55) / SUM(B55:I55)'
$ws.Range('L55').Formula = '=1-K55'

$ws.Range('A56').Value = 'Desert'
$ws.Range('B56').Value = 3
$ws.Range('C56').Value = 2
$ws.Range('D56').Value = 0
$ws.Range('E56').Value = 0
$ws.Range('F56').Value = 9
$ws.Range('G56').Value = 1
$ws.Range('H56').Value = 0
$ws.Range('I56').Value = 5
$ws.Range('J56').Formula = '=SUM(B56:I56)/SUM($B$22:$I$29)'
$ws.Range('K56').Formula = '=(SUM(B56:I56) - F56) / SUM(B56:I56)'
$ws.Range('L56').Formula = '=1-K56'

$ws.Range('A57').Value = 'Football Field'
$ws.Range('B57').Value = 4
$ws.Range('C57').Value = 5
$ws.Range('D57').Value = 0
$ws.Range('E57').Value = 0
$ws.Range('F57').Value = 7
$ws.Range('G57').Value = 1
$ws.Range('H57').Value = 0
$ws.Range('I57').Value = 3
$ws.Range('J57').Formula = '=SUM(B57:I57)/SUM($B$22:$I$29)'
$ws.Range('K57').Formula = '=(SUM(B57:I57) - G57) / SUM(B57:I57)'
$ws.Range('L57').Formula = '=1-K57'

$ws.Range('A58').Value = 'Kitchen'
$ws.Range('B58').Value = 1
$ws.Range('C58').Value = 5
$ws.Range('D58').Value = 2
$ws.Range('E58').Value = 0
$ws.Range('F58').Value = 4
$ws.Range('G58').Value = 0
$ws.Range('H58').Value = 8
$ws.Range('I58').Value = 0
$ws.Range('J58').Formula = '=SUM(B58:I58)/SUM($B$22:$I$29)'
$ws.Range('K58').Formula = '=(SUM(B58:I58) - H58) / SUM(B58:I58)'
$ws.Range('L58').Formula = '=1-K58'

$ws.Range('A59').Value = 'Sky'
$ws.Range('B59').Value = 1
$ws.Range('C59').Value = 8
$ws.Range('D59').Value = 0
$ws.Range('E59').Value = 0
$ws.Range('F59').Value = 4
$ws.Range('G59').Value = 0
$ws.Range('H59').Value = 3
$ws.Range('I59').Value = 4
$ws.Range('J59').Formula = '=SUM(B59:I59)/SUM($B$22:$I$29)'
$ws.Range('K59').Formula = '=(SUM(B59:I59) - I59) / SUM(B59:I59)'
$ws.Range('L59').Formula = '=1-K59'

$ws.Range('A60').Value = 'Percentage'
$ws.Range('B60').Formula = '=SUM(B52:B59) / SUM($B$52:$I$59)'
$ws.Range('C60').Formula = '=SUM(C52:C59) / SUM($B$52:$I$59)'
$ws.Range('D60').Formula = '=SUM(D52:D59) / SUM($B$52:$I$59)'
$ws.Range('E60').Formula = '=SUM(E52:E59) / SUM($B$52:$I$59)'
$ws.Range('F60').Formula = '=SUM(F52:F59) / SUM($B$52:$I$59)'
$ws.Range('G60').Formula = '=SUM(G52:G59) / SUM($B$52:$I$59)'
$ws.Range('H60').Formula = '=SUM(H52:H59) / SUM($B$52:$I$59)'
$ws.Range('I60').Formula = '=SUM(I52:I59) / SUM($B$52:$I$59)'

$ws.Range('A61').Value = 'False Pos'
$ws.Range('B61').Formula = '=(SUM(B52:B59) - B52) / SUM(B52:B59)'
$ws.Range('C61').Formula = '=(SUM(C52:C59) - C53) / SUM(C52:C59)'
$ws.Range('D61').Formula = '=(SUM(D52:D59) - D54) / SUM(D52:D59)'
$ws.Range('F61').Formula = '=(SUM(F52:F59) - F56) / SUM(F52:F59)'
$ws.Range('G61').Formula = '=(SUM(G52:G59) - G57) / SUM(G52:G59)'
$ws.Range('H61').Formula = '=(SUM(H52:H59) - H58) / SUM(H52:H59)'
$ws.Range('I61').Formula = '=(SUM(I52:I59) - I59) / SUM(I52:I59)'
$ws.Range('K61').Value = 'Accuracy'
$ws.Range('K61').Font.Bold = $true
$ws.Range('L61').Formula = '=(B52+C53+D54+E55+F56+G57+H58+I59) / SUM(B52:I59)'

$ws.Range('A62').Value = 'True Neg'
$ws.Range('B62').Formula = '=1-B61'
$ws.Range('C62').Formula = '=1-C61'
$ws.Range('D62').Formula = '=1-D61'
$ws.Range('F62').Formula = '=1-F61'
$ws.Range('G62').Formula = '=1-G61'
$ws.Range('H62').Formula = '=1-H61'
$ws.Range('I62').Formula = '=1-I61'

# ---- Row 65: Baseline label (moved down from old row 35) ----
$ws.Range('A65').Value = 'Baseline Multiresolution Histograms'
$ws.Range('A65').Font.Bold = $true

# ---- Update view/selection state ----
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range('I57').Select() | Out-Null

